$d = $word.ActiveDocument

# --- 1) "2.7) Draw" paragraph: split " Draw" into " " + "Draw" (strike the word "Draw") ---
$pDraw = $d.Paragraphs.Item(73)
$drawText = $pDraw.Range.Text.TrimEnd()
if ($drawText -ne "2.7) Draw") {
    throw "Unexpected paragraph 73 text: [$drawText]"
}
$pDrawStart = $pDraw.Range.Start
$drawWordStart = $pDrawStart + 5   # "2.7) " is 5 chars -> "Draw" begins here
$drawWordEnd = $pDrawStart + 9     # "Draw" is 4 chars
$rDrawWord = $d.Range($drawWordStart, $drawWordEnd)
if ($rDrawWord.Text -ne "Draw") {
    throw "Unexpected Draw sub-range text: [$($rDrawWord.Text)]"
}
$rDrawWord.Font.StrikeThrough = 1

# --- Remove the old (hidden) _GoBack bookmark from its current position ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2) "Loging moves" paragraph: strike the whole paragraph (runs + paragraph mark) ---
$pLoging = $d.Paragraphs.Item(79)
$logingText = $pLoging.Range.Text.TrimEnd()
if ($logingText -ne "Loging moves") {
    throw "Unexpected paragraph 79 text: [$logingText]"
}
$pLoging.Range.Font.StrikeThrough = 1

# --- 3) "MVVM for pages" paragraph: re-insert the _GoBack bookmark right after the run ---
$pMvvm = $d.Paragraphs.Item(80)
$mvvmText = $pMvvm.Range.Text.TrimEnd()
if ($mvvmText -ne "MVVM for pages") {
    throw "Unexpected paragraph 80 text: [$mvvmText]"
}
$mvvmEnd = $pMvvm.Range.End - 1   # position right after "pages", before the paragraph mark

# Placing a bookmark directly at a collapsed (start==end) range sitting exactly on a
# paragraph-end boundary is mishandled by this host, so nudge around it: insert a
# throwaway character at that position, anchor the bookmark there, then remove the
# throwaway character again. The zero-width bookmark stays put.
$insertion = $d.Range($mvvmEnd, $mvvmEnd)
$insertion.InsertAfter("X")

$bmRange = $d.Range($mvvmEnd, $mvvmEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

$markerRange = $d.Range($mvvmEnd, $mvvmEnd + 1)
if ($markerRange.Text -ne "X") {
    throw "Unexpected marker text: [$($markerRange.Text)]"
}
$markerRange.Delete()

Write-Output "done"
